$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot values.
# Numeric-looking Price (column D) values are entered with a leading apostrophe
# (Excel's text quote-prefix) so they are stored as text, matching the sheet's
# existing inline-string formatting instead of being auto-converted to numbers.

$ws.Range("D2").Value = "26.771.27"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.725.85"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'0.9975"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'240.45"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'0.9982"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.4815"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "'0.2591"
$ws.Range("D9").Value = "'0.06174"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "1.723.24"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'15.84"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "'0.06842"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "'0.6030"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "'4.459"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "'76.85"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "'0.9983"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.578.35"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.9979"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'0.000007124"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "1.945.76"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'4.409"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "'8.481"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'5.057"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "'139.86"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D27").Value = "'1.780"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").Value = "'106.38"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'1.368"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").Value = "'3.982"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "'3.667"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'2.589"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "'0.9990"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'0.6172"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'0.9279"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "'2.443"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").Value = "'1.990"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").Value = "'0.9977"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "'5.605"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").Value = "'99.74"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'0.3828"
$ws.Range("D45").Value = "'6.768"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'7.931"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'30.08"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "'51.49"
$ws.Range("E51").Value = "  +0.90%  "
